# Daily attendance processing - 2025-10-13 13:52:24
# Reverse the order of names in the "Recorded By" (column G) list for each row
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text
    if ($text -and $text.Contains(",")) {
        $parts = $text -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newText = [string]::Join(", ", $reversed)
            if ($newText -ne $text) {
                $cell.Value = $newText
            }
        }
    }
}
